$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53 corresponds to 25.9.18 (tuesday) -> Task column (J) was blank, now "off"
$ws.Range("J53").Value = "off"

# Row 54 corresponds to 26.9.18 (wednesday) -> Time In / Time Out / Task now filled in
$ws.Range("G54").Value = 0.4375
$ws.Range("H54").Value = 0.770833333333333
$ws.Range("J54").Value = "watched tutorials of using SQL with oracle SQL developer and loops and function in php"

# Keep the selection/view consistent with the edited cell
$ws.Range("H55").Select
